$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns before D (new quarters), shifting existing data to F:M
$ws.Columns("D:E").Insert()

# 2) Copy number format (incl. font) from column F into the two new columns D:E for every row
#    so the inserted cells pick up the same date/number style as their row, not column C's style.
for ($r = 5; $r -le 102; $r++) {
  $src = $ws.Cells.Item($r, 6)
  $dst = $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 5))
  $src.Copy()
  $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 3) Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with reported figures
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 188600
$ws.Cells.Item(8,5).Value = 173800
$ws.Cells.Item(9,4).Value = 44100
$ws.Cells.Item(9,5).Value = 39300
$ws.Cells.Item(10,4).Value = 144500
$ws.Cells.Item(10,5).Value = 134500
$ws.Cells.Item(12,4).Value = 27200
$ws.Cells.Item(12,5).Value = 26300
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(17,4).Value = 192000
$ws.Cells.Item(17,5).Value = 180900
$ws.Cells.Item(18,4).Value = -3400
$ws.Cells.Item(18,5).Value = -7100
$ws.Cells.Item(20,4).Value = 2600
$ws.Cells.Item(20,5).Value = 2600
$ws.Cells.Item(21,4).Value = 5200
$ws.Cells.Item(21,5).Value = 1200
$ws.Cells.Item(22,4).Value = 4900
$ws.Cells.Item(22,5).Value = 4900
$ws.Cells.Item(23,4).Value = -5800
$ws.Cells.Item(23,5).Value = -9400
$ws.Cells.Item(24,4).Value = -100
$ws.Cells.Item(24,5).Value = 100
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = -5700
$ws.Cells.Item(26,5).Value = -9500
$ws.Cells.Item(27,4).Value = -5700
$ws.Cells.Item(27,5).Value = -9500
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = "NA"
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = -2600
$ws.Cells.Item(32,5).Value = -2600
$ws.Cells.Item(33,4).Value = -5700
$ws.Cells.Item(33,5).Value = -9500
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = -5700
$ws.Cells.Item(35,5).Value = -9500
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 566300
$ws.Cells.Item(41,5).Value = 577300
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,4).Value = 94400
$ws.Cells.Item(43,5).Value = 74200
$ws.Cells.Item(44,4).Value = 200
$ws.Cells.Item(44,5).Value = 200
$ws.Cells.Item(45,4).Value = 46600
$ws.Cells.Item(45,5).Value = 50500
$ws.Cells.Item(46,4).Value = 707500
$ws.Cells.Item(46,5).Value = 702200
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(48,4).Value = 70200
$ws.Cells.Item(48,5).Value = 60200
$ws.Cells.Item(49,4).Value = 50700
$ws.Cells.Item(49,5).Value = 17800
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 65900
$ws.Cells.Item(52,5).Value = 61700
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 894300
$ws.Cells.Item(54,5).Value = 841800
$ws.Cells.Item(57,4).Value = 10100
$ws.Cells.Item(57,5).Value = 6000
$ws.Cells.Item(58,4).Value = "NA"
$ws.Cells.Item(58,5).Value = 900
$ws.Cells.Item(59,4).Value = 189200
$ws.Cells.Item(59,5).Value = 165700
$ws.Cells.Item(60,4).Value = 199400
$ws.Cells.Item(60,5).Value = 172700
$ws.Cells.Item(61,4).Value = 366600
$ws.Cells.Item(61,5).Value = 364500
$ws.Cells.Item(62,4).Value = 10800
$ws.Cells.Item(62,5).Value = 5200
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 576700
$ws.Cells.Item(66,5).Value = 542400
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = -235700
$ws.Cells.Item(72,5).Value = -230000
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 317600
$ws.Cells.Item(76,5).Value = 299400
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = -5700
$ws.Cells.Item(81,5).Value = -9500
$ws.Cells.Item(83,4).Value = 6100
$ws.Cells.Item(83,5).Value = 5700
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 23000
$ws.Cells.Item(89,5).Value = 20200
$ws.Cells.Item(91,4).Value = -3300
$ws.Cells.Item(91,5).Value = -2800
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -39000
$ws.Cells.Item(94,5).Value = -8700
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = 5300
$ws.Cells.Item(100,5).Value = -600
$ws.Cells.Item(101,4).Value = -200
$ws.Cells.Item(101,5).Value = -900
$ws.Cells.Item(102,4).Value = -11000
$ws.Cells.Item(102,5).Value = 10000

# 4) A handful of prior-quarter figures were restated in this filing; update those cells
#    (now shifted two columns to the right, i.e. old F/G -> new H/I, etc.)
$ws.Cells.Item(8,8).Value = 141200
$ws.Cells.Item(8,9).Value = 130300
$ws.Cells.Item(10,8).Value = 107600
$ws.Cells.Item(10,9).Value = 99500
$ws.Cells.Item(17,8).Value = 141400
$ws.Cells.Item(17,9).Value = 130500
$ws.Cells.Item(18,8).Value = -200
$ws.Cells.Item(18,9).Value = -200
$ws.Cells.Item(21,8).Value = 4200
$ws.Cells.Item(21,9).Value = 4400
$ws.Cells.Item(23,8).Value = 0
$ws.Cells.Item(23,9).Value = 300
$ws.Cells.Item(24,8).Value = -33200
$ws.Cells.Item(26,8).Value = 33100
$ws.Cells.Item(26,9).Value = 300
$ws.Cells.Item(27,8).Value = 33100
$ws.Cells.Item(27,9).Value = 300
$ws.Cells.Item(29,8).Value = -33300
$ws.Cells.Item(33,8).Value = -100
$ws.Cells.Item(33,9).Value = 300
$ws.Cells.Item(35,8).Value = -100
$ws.Cells.Item(35,9).Value = 300
$ws.Cells.Item(81,8).Value = -100
$ws.Cells.Item(81,9).Value = 300
$ws.Cells.Item(91,6).Value = -2600
$ws.Cells.Item(91,7).Value = -2800
$ws.Cells.Item(91,8).Value = -2000
$ws.Cells.Item(91,9).Value = -1900
$ws.Cells.Item(91,10).Value = -1800
